# Split the single run "Fitting options:" into three runs, changing
# "options" to "approach" along the way:
#   <w:r><w:t>Fitting options:</w:t></w:r>
# becomes
#   <w:r><w:t xml:space="preserve">Fitting </w:t></w:r>
#   <w:r><w:t>approach</w:t></w:r>
#   <w:r><w:t>:</w:t></w:r>

$d = $word.ActiveDocument

$find = $d.Content
$find.Find.Execute("Fitting options:", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)

$start = $find.Start
$end = $find.End

# Re-seat a fresh Range over the matched text (rather than reusing the Find's
# own Range object) and feed it the replacement as three sibling runs via
# InsertXML -- this swaps the whole match for the new run sequence in one
# shot, without going through a separate Delete() call (which would cause
# the adjacent, identically-formatted runs to be silently coalesced back
# into a single run).
$target = $d.Range($start, $end)

$packageXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Fitting </w:t></w:r>
            <w:r><w:t>approach</w:t></w:r>
            <w:r><w:t>:</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($packageXml)
